$d = $word.ActiveDocument

# Locate the paragraph that ends the "Samenwerken/Overleg/Afspraken" bullet
# block ("Wat vindt de rest?") - the new "Analyse" / "Test verslag" items
# are inserted right after it, before the page break.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Wat vindt de rest?*") {
        $anchor = $p
        break
    }
}

$items = @(
    @{ Text = "Analyse"; Level = 1 },
    @{ Text = "SD's maken"; Level = 2 },
    @{ Text = "Klassen diagram bijwerken (Server, Protocol, Client)"; Level = 2 },
    @{ Text = "UCD bijwerken"; Level = 2 },
    @{ Text = "STD maken/bijwerken"; Level = 2 },
    @{ Text = "Test verslag"; Level = 1 }
)

foreach ($item in $items) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $item.Text
    $anchor.Range.ListFormat.ListLevelNumber = $item.Level
}
